$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 to reuse the same password value as B2 ("Test@1234")
$ws.Range("B3").Value = $ws.Range("B2").Value2

# Move the active selection to L10
$ws.Range("L10").Select()
